$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The only real content change: E8 goes from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update selection to match the post-edit state (cursor left on E8)
$ws.Range("E8").Select()
